# Update latest output (run 93)
# - "Schedule" sheet: row 2 gets new aggregated values, row 3 is removed.
# - "Detailed" sheet: a new first interval is inserted at row 2 and all
#   subsequent rows shift down by one (Price/Type/Pump_Status refreshed
#   with the new run's values), with one new row appended at the end.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Schedule" (first sheet)
# -----------------------------------------------------------------
$sched = $wb.Worksheets.Item(1)

$sched.Cells.Item(2, 1).Value = 46041                  # A2 Start Time
$sched.Cells.Item(2, 3).Value = 16                     # C2 Duration (h)
$sched.Cells.Item(2, 4).Value = 60.48                  # D2 Volume (ML)
$sched.Cells.Item(2, 5).Value = 723.6466672500002      # E2 Cost ($)
$sched.Cells.Item(2, 6).Value = 11.96505732886905      # F2 Unit Cost ($/ML)

# Row 3 no longer exists in the refreshed output
$sched.Rows.Item(3).Delete()

# -----------------------------------------------------------------
# Sheet "Detailed" (second sheet)
# -----------------------------------------------------------------
$det = $wb.Worksheets.Item(2)

# Full refreshed dataset for rows 2..49 (DateTime, Price, Type, Date, Pump_Status)
$data = @(
  @(46041,                 "57.06003",            "historical", 46041, "ON"),
  @(46041.02083333334,     "57.06003",            "historical", 46041, "ON"),
  @(46041.04166666666,     "56.98",               "historical", 46041, "ON"),
  @(46041.0625,            "48.2193",             "historical", 46041, "ON"),
  @(46041.08333333334,     "35.87995",            "historical", 46041, "ON"),
  @(46041.10416666666,     "36.2",                "forecast",   46041, "ON"),
  @(46041.125,             "48.63149",            "forecast",   46041, "ON"),
  @(46041.14583333334,     "56.98",               "forecast",   46041, "ON"),
  @(46041.16666666666,     "57.06003",            "forecast",   46041, "ON"),
  @(46041.1875,            "57.08",               "forecast",   46041, "ON"),
  @(46041.20833333334,     "59.24143",            "forecast",   46041, "ON"),
  @(46041.22916666666,     "61.17605",            "forecast",   46041, "ON"),
  @(46041.25,              "60.35837",            "forecast",   46041, "ON"),
  @(46041.27083333334,     "57.06003",            "forecast",   46041, "ON"),
  @(46041.29166666666,     "32.5261",             "forecast",   46041, "ON"),
  @(46041.3125,            "1.16848",             "forecast",   46041, "ON"),
  @(46041.33333333334,     "0.51",                "forecast",   46041, "ON"),
  @(46041.35416666666,     "6.82594",             "forecast",   46041, "ON"),
  @(46041.375,             "0.51",                "forecast",   46041, "ON"),
  @(46041.39583333334,     "-4.91542",            "forecast",   46041, "ON"),
  @(46041.41666666666,     "-5.50985",            "forecast",   46041, "ON"),
  @(46041.4375,            "-5.74313",            "forecast",   46041, "ON"),
  @(46041.45833333334,     "-5.66734",            "forecast",   46041, "ON"),
  @(46041.47916666666,     "-5.50985",            "forecast",   46041, "ON"),
  @(46041.5,               "-5.74313",            "forecast",   46041, "ON"),
  @(46041.52083333334,     "-5.78385",            "forecast",   46041, "ON"),
  @(46041.54166666666,     "-1.39991",            "forecast",   46041, "ON"),
  @(46041.5625,            "-5.01",               "forecast",   46041, "ON"),
  @(46041.58333333334,     "-2.63472",            "forecast",   46041, "ON"),
  @(46041.60416666666,     "-0.9186299999999999", "forecast",   46041, "ON"),
  @(46041.625,             "0.00031",             "forecast",   46041, "ON"),
  @(46041.64583333334,     "0.51",                "forecast",   46041, "ON"),
  @(46041.66666666666,     "0.7",                 "forecast",   46041, "OFF"),
  @(46041.6875,            "-2.47487",            "forecast",   46041, "OFF"),
  @(46041.70833333334,     "-2.56535",            "forecast",   46041, "OFF"),
  @(46041.72916666666,     "0.00958",             "forecast",   46041, "OFF"),
  @(46041.75,              "12.00955",            "forecast",   46041, "OFF"),
  @(46041.77083333334,     "43.32859",            "forecast",   46041, "OFF"),
  @(46041.79166666666,     "56.78845",            "forecast",   46041, "OFF"),
  @(46041.8125,            "60.73383",            "forecast",   46041, "OFF"),
  @(46041.83333333334,     "64.10924",            "forecast",   46041, "OFF"),
  @(46041.85416666666,     "65",                  "forecast",   46041, "OFF"),
  @(46041.875,             "61.80428",            "forecast",   46041, "OFF"),
  @(46041.89583333334,     "62.51999",            "forecast",   46041, "OFF"),
  @(46041.91666666666,     "61.16311",            "forecast",   46041, "OFF"),
  @(46041.9375,            "58.63955",            "forecast",   46041, "OFF"),
  @(46041.95833333334,     "57.62834",            "forecast",   46041, "OFF"),
  @(46041.97916666666,     "61.64359",            "forecast",   46041, "OFF")
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 2
  $row = $data[$i]
  $det.Cells.Item($r, 1).Value = $row[0]
  $det.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
  $det.Cells.Item($r, 2).Value = $row[1]
  $det.Cells.Item($r, 3).Value = $row[2]
  $det.Cells.Item($r, 4).Value = $row[3]
  $det.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD"
  $det.Cells.Item($r, 5).Value = $row[4]
}
